$wb = $excel.ActiveWorkbook

# --- Sheet R1: update Elapsed Duration(Hrs) values ---
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3930:15:41"
$ws.Range("G3").Value = "69:48:19"
$ws.Range("G4").Value = "92:48:19"

# --- Sheet R2: update Elapsed Duration(Hrs) values ---
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12111:39:21"
$ws.Range("G3").Value = "3241:22:50"
$ws.Range("G4").Value = "479:34:24"

# --- Sheet R3: remove row 3 (HAJ0125 / Haj Removal) entirely ---
$ws = $wb.Worksheets.Item("R3")
$ws.Rows.Item(3).Delete()

# --- Sheet R4: update Elapsed Duration(Hrs) values ---
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2957:29:10"
$ws.Range("G3").Value = "184:41:25"
$ws.Range("G4").Value = "72:53:50"
$ws.Range("G5").Value = "70:31:23"

# --- Sheet R5: update Elapsed Duration(Hrs) value and remove row 3 (MAK0875) ---
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "431:28:09"
$ws.Rows.Item(3).Delete()

# --- Sheet R6: update Elapsed Duration(Hrs) value ---
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "72:00:27"
